$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# Append three more days of tracking data (2021-03-09 .. 2021-03-11),
# continuing the existing daily log with all-zero measurements.
# Copy the last data row (49) down into rows 50:52 first so the new rows
# inherit the same date formatting (style) as the rest of column A, then
# overwrite the copied values with the correct dates and zeros.
$ws.Rows("49").Copy()
$ws.Rows("50:52").Insert()

$dates = @(44264, 44265, 44266)
for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = 50 + $i
    $ws.Cells.Item($row, 1).Value2 = $dates[$i]
    for ($col = 2; $col -le 7; $col++) {
        $ws.Cells.Item($row, $col).Value2 = 0
    }
}

$excel.CutCopyMode = $false

# Match the final selection/scroll state left in the saved workbook.
$ws.Range("I46").Select()
